$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Mark "Uzupelnienie tabeli z emailami" (row 14), "Usuwanie filmow" (row 16)
# and "Usuwanie muzyki" (row 18) as finished ("Zakonczone") instead of
# "Nie rozpoczete" (not started).
$ws.Range("F14").Value = "Zakończone"
$ws.Range("F16").Value = "Zakończone"
$ws.Range("F18").Value = "Zakończone"

# Update the active selection on the sheet.
$ws.Range("H16").Select()
